# Updated symbol list on Mon Dec 12 21:47:18 UTC 2022 with GitHub Actions
#
# Applies the per-cell value updates to the "cryptos" price table on
# Sheet1. All target cells are plain text cells (t="inlineStr" in the
# original OOXML) holding numeric-looking strings (prices) or plain
# strings (coin names / links / rank labels), so each write forces the
# cell to Text format before assigning the value and then restores the
# "Normal" style afterwards so no stray number-format / style index is
# left behind on cells that originally had none.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        $range,
        [string]$value
    )
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "276.35"
Set-TextValue $ws.Range("D3") "21.14"
Set-TextValue $ws.Range("D4") "6.260"
Set-TextValue $ws.Range("D5") "0.06240"
Set-TextValue $ws.Range("D6") "3.548"
Set-TextValue $ws.Range("D7") "1.532"
Set-TextValue $ws.Range("D8") "6.563"
Set-TextValue $ws.Range("D9") "0.8249"
Set-TextValue $ws.Range("D10") "0.1657"
Set-TextValue $ws.Range("D11") "0.08274"
Set-TextValue $ws.Range("D12") "0.03509"
Set-TextValue $ws.Range("D13") "0.03161"
Set-TextValue $ws.Range("D14") "0.09141"
Set-TextValue $ws.Range("D15") "3.759"
Set-TextValue $ws.Range("D16") "0.001643"
Set-TextValue $ws.Range("D17") "0.04683"
Set-TextValue $ws.Range("D18") "0.006239"
Set-TextValue $ws.Range("D19") "0.006225"
Set-TextValue $ws.Range("D22") "3.726"
Set-TextValue $ws.Range("D24") "0.01398"
Set-TextValue $ws.Range("D28") "0.0002730"
Set-TextValue $ws.Range("D40") "0.04742"
Set-TextValue $ws.Range("B41") "KickToken"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws.Range("D41") "0.007019"
Set-TextValue $ws.Range("E41") "40KickTokenKICK"
Set-TextValue $ws.Range("B42") "BKEXToken"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Range("D42") "0.1123"
Set-TextValue $ws.Range("E42") "41BKEXTokenBKK"
Set-TextValue $ws.Range("B43") "CEJI"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws.Range("D43") "0.003514"
Set-TextValue $ws.Range("E43") "42CEJICEJI"
Set-TextValue $ws.Range("D45") "0.00006246"
Set-TextValue $ws.Range("D46") "0.00000000748"
Set-TextValue $ws.Range("D47") "0.7212"
Set-TextValue $ws.Range("D48") "0.001393"
Set-TextValue $ws.Range("D49") "0.00001895"
Set-TextValue $ws.Range("E49") "48CryptobidCoinCBCBestin24h"
Set-TextValue $ws.Range("D50") "0.01237"
